$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 18 and row 19 and must be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $addr18 = "$col" + "18"
    $addr19 = "$col" + "19"
    $v18 = $ws.Range($addr18).Value()
    $v19 = $ws.Range($addr19).Value()
    $ws.Range($addr18).Value = $v19
    $ws.Range($addr19).Value = $v18
}
